# Apply the crypto price/volume refresh described by the commit diff.
# Updates columns B (Coin), C (Link), D (Price) and E (Volume(1h)) for rows 2-51
# of the active worksheet to the new scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as literal text (e.g. "1.00", "0.670",
# "76.255.17") rather than numbers, so that trailing zeros, thousands-style
# separators and multi-dot values are preserved exactly. For the new prices
# that would otherwise be auto-recognised as a plain number, force the cell
# to Text format *before* writing the value so Excel keeps it as a string.
$textPriceCells = @(
    "D4", "D5", "D6", "D8", "D11", "D13", "D16", "D20", "D21", "D22", "D23", "D24", "D27",
    "D30", "D32", "D35", "D36", "D37", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47",
    "D48", "D50", "D51"
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# New cell values scraped for this update.

# Row 2
$ws.Range("D2").Value = "76.255.17"
$ws.Range("E2").Value = "  +1.90%  "

# Row 3
$ws.Range("D3").Value = "2.856.25"
$ws.Range("E3").Value = "  +7.26%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "195.60"
$ws.Range("E5").Value = "  +5.07%  "

# Row 6
$ws.Range("D6").Value = "600.72"
$ws.Range("E6").Value = "  +2.66%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +3.72%  "

# Row 9
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
$ws.Range("D10").Value = "2.850.14"
$ws.Range("E10").Value = "  +6.99%  "

# Row 11
$ws.Range("D11").Value = "0.392"
$ws.Range("E11").Value = "  +10.60%  "

# Row 12
$ws.Range("E12").Value = "  -1.96%  "

# Row 13
$ws.Range("D13").Value = "4.90"
$ws.Range("E13").Value = "  +3.68%  "

# Row 14
$ws.Range("D14").Value = "3.377.78"
$ws.Range("E14").Value = "  +7.09%  "

# Row 15
$ws.Range("D15").Value = "76.067.12"
$ws.Range("E15").Value = "  +1.86%  "

# Row 16
$ws.Range("D16").Value = "27.52"
$ws.Range("E16").Value = "  +4.15%  "

# Row 17
$ws.Range("E17").Value = "  +1.88%  "

# Row 18
$ws.Range("D18").Value = "2.848.60"
$ws.Range("E18").Value = "  +7.05%  "

# Row 19
$ws.Range("E19").Value = "  -0.90%  "

# Row 20
$ws.Range("D20").Value = "12.43"
$ws.Range("E20").Value = "  +4.83%  "

# Row 21
$ws.Range("D21").Value = "382.73"
$ws.Range("E21").Value = "  +3.22%  "

# Row 22
$ws.Range("D22").Value = "2.35"
$ws.Range("E22").Value = "  +4.35%  "

# Row 23
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  +1.52%  "

# Row 24
$ws.Range("D24").Value = "71.91"
$ws.Range("E24").Value = "  +3.55%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("D26").Value = "3.008.41"
$ws.Range("E26").Value = "  +7.45%  "

# Row 27
$ws.Range("D27").Value = "4.22"
$ws.Range("E27").Value = "  +2.28%  "

# Row 28
$ws.Range("E28").Value = "  +4.32%  "

# Row 29
$ws.Range("E29").Value = "  +11.38%  "

# Row 30
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("E31").Value = "  +0.30%  "

# Row 32
$ws.Range("D32").Value = "513.26"
$ws.Range("E32").Value = "  -1.17%  "

# Row 33
$ws.Range("E33").Value = "  +1.25%  "

# Row 34
$ws.Range("E34").Value = "  +4.35%  "

# Row 35
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36
$ws.Range("D36").Value = "167.34"
$ws.Range("E36").Value = "  +2.93%  "

# Row 37
$ws.Range("D37").Value = "20.02"
$ws.Range("E37").Value = "  +4.43%  "

# Row 38
$ws.Range("E38").Value = "  +0.48%  "

# Row 39
$ws.Range("D39").Value = "19.48"
$ws.Range("E39").Value = "  +0.77%  "

# Row 40
$ws.Range("D40").Value = "186.03"
$ws.Range("E40").Value = "  +9.88%  "

# Row 41
$ws.Range("E41").Value = "  -0.10%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.09"
$ws.Range("E42").Value = "  +2.63%  "

# Row 43
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.344"
$ws.Range("E43").Value = "  +4.83%  "

# Row 44
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  +0.88%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.0912"
$ws.Range("E45").Value = "  +8.56%  "

# Row 46
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "1.23"
$ws.Range("E46").Value = "  +4.49%  "

# Row 47
$ws.Range("D47").Value = "40.24"
$ws.Range("E47").Value = "  +3.13%  "

# Row 48
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +1.95%  "

# Row 49
$ws.Range("E49").Value = "  +9.67%  "

# Row 50
$ws.Range("D50").Value = "0.670"
$ws.Range("E50").Value = "  +13.59%  "

# Row 51
$ws.Range("D51").Value = "3.74"
$ws.Range("E51").Value = "  +3.13%  "

Write-Host "Updated $($textPriceCells.Count) text-formatted price cells and $($ws.UsedRange.Rows.Count) rows total."
